$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-25 13:10:26"
$wsZhCn.Range("G2").Value = "2016-01-25 13:11:11"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-25 13:10:36"
$wsDeDe.Range("G2").Value = "2016-01-25 13:11:29"
